$wb = $excel.ActiveWorkbook

# --- update selection on the existing "menuData" sheet (it stops being the active tab) ---
$menuData = $wb.Worksheets.Item("menuData")
$menuData.Range("G9").Select()

# --- add "orderdata" sheet (after menuData) ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$orderdata = $wb.Worksheets.Add($null, $after)
$orderdata.Name = "orderdata"
$orderdata.Range("A1").Value = "menuType"
$orderdata.Range("B1").Value = "itemName"
$orderdata.Range("A2").Value = "HBO Cafe Grill"
$orderdata.Range("B2").Value = "Bacon"
$orderdata.Range("B2").Select()

# --- add "instructions" sheet ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$instructions = $wb.Worksheets.Add($null, $after)
$instructions.Name = "instructions"
$instructions.Range("A1").Value = "instructions"
$instructions.Range("A2").Value = "Placing order for 2 items from non veg"
$instructions.Range("A3").Select()

# --- add "cardDetails" sheet ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$cardDetails = $wb.Worksheets.Add($null, $after)
$cardDetails.Name = "cardDetails"
$cardDetails.Range("A2").Value = "'5399999999999999"
$cardDetails.Range("A1").Value = "cardnumber"
$cardDetails.Range("B1").Value = "mm"
$cardDetails.Range("C1").Value = "year"
$cardDetails.Range("D1").Value = "cvv"
$cardDetails.Range("E1").Value = "zipcode"
$cardDetails.Range("B2").Value = 11
$cardDetails.Range("C2").Value = 2013
$cardDetails.Range("D2").Value = 999
$cardDetails.Range("E2").Value = "'98760"
$cardDetails.Range("F8").Select()

# --- add "review" sheet ---
$after = $wb.Worksheets.Item($wb.Worksheets.Count)
$review = $wb.Worksheets.Add($null, $after)
$review.Name = "review"
$review.Range("A2").Value = "Good menu"
$review.Range("A1").Value = "reviewcomments"
$review.Range("B1").Value = "menuorder"
$review.Range("B2").Value = "HBO"
$review.Range("D7").Select()
